# Reposition the six data-callout text boxes ("tx7".."tx12") that live
# inside the single diagram group on slide 1. Only their x/y offsets
# change; widths/heights are untouched.
#
# Shape.Left/.Top are expressed in points while the XML stores EMU
# (1 pt = 12700 EMU), and PowerPoint's COM layer keeps Left/Top as
# single-precision floats. Converting EMU -> points -> EMU can
# therefore truncate to one EMU below the intended value, so a tiny
# half-EMU epsilon is added before the conversion to land back on the
# exact target EMU after the round trip.

$EmuPerPt = 12700
$halfEmuInPt = 0.5 / $EmuPerPt

function ConvertTo-Points([double]$emu) {
    return ($emu / $EmuPerPt) + $halfEmuInPt
}

$newPositionsEmu = @{
    "tx7"  = @{ x = 5855948; y = 4583028 }
    "tx8"  = @{ x = 6151045; y = 4889412 }
    "tx9"  = @{ x = 3717972; y = 3490613 }
    "tx10" = @{ x = 4025137; y = 3840083 }
    "tx11" = @{ x = 4068564; y = 3076399 }
    "tx12" = @{ x = 4158867; y = 3423329 }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(2)

foreach ($item in $g.GroupItems) {
    $pos = $newPositionsEmu[$item.Name]
    if ($pos) {
        $item.Left = ConvertTo-Points $pos.x
        $item.Top  = ConvertTo-Points $pos.y
    }
}
